$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 18
$ws.Range("H18").Value = 1225.2667
$ws.Range("I18").Value = 1038.9
$ws.Range("J18").Value = 1598
$ws.Range("K18").Value = 1038.9
$ws.Range("L18").Value = 1598
$ws.Range("M18").Value = -754.9000000000001
$ws.Range("N18").Value = -2166

# ALC row 40
$ws.Range("H40").Value = 4946.2812
$ws.Range("I40").Value = 7617.4375
$ws.Range("J40").Value = 2275.125
$ws.Range("K40").Value = 7617.4375
$ws.Range("L40").Value = 2275.125
$ws.Range("M40").Value = -7442.4375
$ws.Range("N40").Value = -2625.125

# ALC row 74
$ws.Range("H74").Value = 6197.6875
$ws.Range("J74").Value = 4920
$ws.Range("L74").Value = 4920
$ws.Range("N74").Value = -6792

# ALC row 77
$ws.Range("H77").Value = 6197.6875
$ws.Range("J77").Value = 4920
$ws.Range("L77").Value = 24600
$ws.Range("N77").Value = -33960

# ALC row 129
$ws.Range("H129").Value = 2573
$ws.Range("I129").Value = 359.4
$ws.Range("J129").Value = 3363.5715
$ws.Range("K129").Value = 1078.2
$ws.Range("L129").Value = 10090.7145
$ws.Range("M129").Value = 3921.8
$ws.Range("N129").Value = -20090.7145

# ALC row 132
$ws.Range("H132").Value = 5210226
$ws.Range("I132").Value = 10417466
$ws.Range("J132").Value = 2986.6667
$ws.Range("K132").Value = 31252398
$ws.Range("L132").Value = 8960.000100000001
$ws.Range("M132").Value = -31249868
$ws.Range("N132").Value = -14020.0001

# ALC row 137
$ws.Range("H137").Value = 998.8182
$ws.Range("I137").Value = 855.1429
$ws.Range("J137").Value = 1250.25
$ws.Range("K137").Value = 2565.4287
$ws.Range("L137").Value = 3750.75
$ws.Range("M137").Value = -15.42870000000039
$ws.Range("N137").Value = -8850.75

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 5687.7407
$ws.Range("I32").Value = 2922.125
$ws.Range("J32").Value = 27812.666
$ws.Range("K32").Value = 2922.125
$ws.Range("L32").Value = 27812.666
$ws.Range("M32").Value = -2635.125
$ws.Range("N32").Value = -28386.666

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Range("H105").Value = 1264452.9
$ws.Range("I105").Value = 1749911.8
$ws.Range("K105").Value = 1749911.8
$ws.Range("M105").Value = -1748164.8

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Range("H7").Value = 89.92857
$ws.Range("I7").Value = 99.888885
$ws.Range("J7").Value = 72
$ws.Range("K7").Value = 99.888885
$ws.Range("L7").Value = 72
$ws.Range("M7").Value = 13.111115
$ws.Range("N7").Value = -298

# CRP row 17
$ws.Range("H17").Value = 14260
$ws.Range("I17").Value = 3900
$ws.Range("J17").Value = 21166.666
$ws.Range("K17").Value = 3900
$ws.Range("L17").Value = 21166.666
$ws.Range("M17").Value = -3726
$ws.Range("N17").Value = -21514.666

# CRP row 22
$ws.Range("H22").Value = 50000348
$ws.Range("I22").Value = 71428770
$ws.Range("J22").Value = 690
$ws.Range("K22").Value = 71428770
$ws.Range("L22").Value = 690
$ws.Range("M22").Value = -71428420
$ws.Range("N22").Value = -1390

# CRP row 25
$ws.Range("H25").Value = 7006.5
$ws.Range("I25").Value = 7000
$ws.Range("K25").Value = 7000
$ws.Range("M25").Value = -6826

# CRP row 41
$ws.Range("H41").Value = 17757.25
$ws.Range("I41").Value = 4609.6665
$ws.Range("J41").Value = 57200
$ws.Range("K41").Value = 4609.6665
$ws.Range("L41").Value = 57200
$ws.Range("M41").Value = -4181.6665
$ws.Range("N41").Value = -58056

# CRP row 50
$ws.Range("H50").Value = 1200
$ws.Range("I50").Value = 1200
$ws.Range("K50").Value = 1200
$ws.Range("M50").Value = -575

# CRP row 51
$ws.Range("H51").Value = 49950
$ws.Range("J51").Value = 49950
$ws.Range("L51").Value = 49950
$ws.Range("N51").Value = -51422

# CRP row 60
$ws.Range("H60").Value = 10596.5625
$ws.Range("J60").Value = 11103
$ws.Range("L60").Value = 11103
$ws.Range("N60").Value = -12125

# CRP row 61
$ws.Range("H61").Value = 49950
$ws.Range("J61").Value = 49950
$ws.Range("L61").Value = 49950
$ws.Range("N61").Value = -50646

# CRP row 74
$ws.Range("H74").Value = 13049.182
$ws.Range("I74").Value = 2285
$ws.Range("J74").Value = 14125.6
$ws.Range("K74").Value = 2285
$ws.Range("L74").Value = 14125.6
$ws.Range("M74").Value = -1411
$ws.Range("N74").Value = -15873.6

# CRP row 77
$ws.Range("H77").Value = 13049.182
$ws.Range("I77").Value = 2285
$ws.Range("J77").Value = 14125.6
$ws.Range("K77").Value = 6855
$ws.Range("L77").Value = 42376.8
$ws.Range("M77").Value = -2487
$ws.Range("N77").Value = -51112.8

# CRP row 105
$ws.Range("H105").Value = 2466.3333
$ws.Range("I105").Value = 1432.6666
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 1432.6666
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = 314.3334
$ws.Range("N105").Value = -6994

# CRP row 134
$ws.Range("H134").Value = 1233.8684
$ws.Range("I134").Value = 1118.7693
$ws.Range("K134").Value = 3356.3079
$ws.Range("M134").Value = -821.3078999999998

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4
$ws.Range("H4").Value = 85166.836
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 510001
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 1530003
$ws.Range("M4").Value = -488
$ws.Range("N4").Value = -1530227

# CUL row 107
$ws.Range("H107").Value = 245.68
$ws.Range("I107").Value = 205.66667
$ws.Range("J107").Value = 282.6154
$ws.Range("K107").Value = 617.00001
$ws.Range("L107").Value = 847.8462000000001
$ws.Range("M107").Value = 1302.99999
$ws.Range("N107").Value = -4687.8462

$ws = $wb.Worksheets.Item("LTW")
# LTW row 4
$ws.Range("H4").Value = 5333
$ws.Range("I4").Value = 4000
$ws.Range("J4").Value = 6666
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 6666
$ws.Range("M4").Value = -3887
$ws.Range("N4").Value = -6892

# LTW row 28
$ws.Range("H28").Value = 5333
$ws.Range("I28").Value = 4000
$ws.Range("J28").Value = 6666
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 6666
$ws.Range("M28").Value = -3768
$ws.Range("N28").Value = -7130

# LTW row 37
$ws.Range("H37").Value = 5333
$ws.Range("I37").Value = 4000
$ws.Range("J37").Value = 6666
$ws.Range("K37").Value = 4000
$ws.Range("L37").Value = 6666
$ws.Range("M37").Value = -3893
$ws.Range("N37").Value = -6880

# LTW row 46
$ws.Range("H46").Value = 1188.7778
$ws.Range("I46").Value = 1179.8
$ws.Range("J46").Value = 1200
$ws.Range("K46").Value = 1179.8
$ws.Range("L46").Value = 1200
$ws.Range("M46").Value = -991.8
$ws.Range("N46").Value = -1576

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2
$ws.Range("H2").Value = 1116055.5
$ws.Range("J2").Value = 1116055.5
$ws.Range("L2").Value = 1116055.5
$ws.Range("N2").Value = -1116279.5
